$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.448.33'
$ws.Range("E2").Value = '  +1.99%  '

# Row 3
$ws.Range("D3").Value = '1.671.13'
$ws.Range("E3").Value = '  +1.73%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '219.62'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +2.27%  '

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.5278'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +1.35%  '

# Row 7
$ws.Range("E7").Value = '  -0.05%  '

# Row 8
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.2671'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +2.65%  '

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.06369'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +0.33%  '

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '21.75'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +5.38%  '

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.07795'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +1.61%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.674.70'
$ws.Range("E12").Value = '  +1.87%  '

# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '4.471'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +1.13%  '

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '0.5548'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +1.16%  '

# Row 15
$ws.Range("D15").Value = '0.0₅8282'
$ws.Range("E15").Value = '  +1.22%  '

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '65.44'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +1.46%  '

# Row 17
$ws.Range("D17").Value = '26.453.48'
$ws.Range("E17").Value = '  +1.91%  '

# Row 18
$ws.Range("E18").Value = '  +0.05%  '

# Row 19
$ws.Range("E19").Value = '  +1.12%  '

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '192.72'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +1.85%  '

# Row 21
$ws.Range("E21").Value = '  +1.93%  '

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '6.281'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +0.58%  '

# Row 23
$ws.Range("E23").Value = '  +0.04%  '

# Row 24
$ws.Range("E24").Value = '  +1.74%  '

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '138.41'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -3.60%  '

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '7.392'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +0.36%  '

# Row 27
$ws.Range("E27").Value = '  +2.21%  '

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '1.411'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +0.95%  '

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '0.06171'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +5.02%  '

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '1.289'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +2.17%  '

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '3.611'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +6.59%  '

# Row 32
$ws.Range("E32").Value = '  +0.86%  '

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '1.679'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +2.39%  '

# Row 34
$ws.Range("E34").Value = '  +1.62%  '

# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.6065'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +7.91%  '

# Row 36
$ws.Range("E36").Value = '  +0.96%  '

# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '2.761'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +0.37%  '

# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.01608'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +0.63%  '

# Row 39
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '1.093.72'
$ws.Range("E39").Value = '  +6.62%  '

# Row 40
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '6.024'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +2.98%  '

# Row 41
$ws.Range("E41").Value = '  +0.36%  '

# Row 42
$ws.Range("E42").Value = '  -0.09%  '

# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '100.64'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +2.15%  '

# Row 44
$ws.Range("D44").Value = '1.813.53'
$ws.Range("E44").Value = '  +1.19%  '

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '58.01'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +4.65%  '

# Row 46
$ws.Range("D46").Value = '0.0₈107'
$ws.Range("E46").Value = '  +0.53%  '

# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '8.150'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +1.63%  '

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.9957'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -0.21%  '

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '0.05202'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +1.12%  '

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '1.471'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +7.56%  '

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.4229'
$cell.Style = "Normal"
